$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, centered, bordered) from AW1 into the new header cells AX1:AZ1
$ws.Range("AW1").Copy()
$ws.Range("AX1:AZ1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header dates for the three new quarterly columns
$hdr = New-Object "object[,]" 1,3
$hdr[0,0] = "31/12/2023"
$hdr[0,1] = "31/03/2024"
$hdr[0,2] = "30/06/2024"
$ws.Range("AX1:AZ1").Value = $hdr

$row2 = New-Object "object[,]" 1,3
$row2[0,0] = 349180.992
$row2[0,1] = 361632
$row2[0,2] = 357748
$ws.Range("AX2:AZ2").Value = $row2
$row3 = New-Object "object[,]" 1,3
$row3[0,0] = 122833
$row3[0,1] = 131762
$row3[0,2] = 105762
$ws.Range("AX3:AZ3").Value = $row3
$row4 = New-Object "object[,]" 1,3
$row4[0,0] = 26100
$row4[0,1] = 32190
$row4[0,2] = 37949
$ws.Range("AX4:AZ4").Value = $row4
$row5 = New-Object "object[,]" 1,3
$row5[0,0] = 33478
$row5[0,1] = 34386
$row5[0,2] = 0
$ws.Range("AX5:AZ5").Value = $row5
$row6 = New-Object "object[,]" 1,3
$row6[0,0] = 20163
$row6[0,1] = 24335
$row6[0,2] = 26080
$ws.Range("AX6:AZ6").Value = $row6
$row7 = New-Object "object[,]" 1,3
$row7[0,0] = 27640
$row7[0,1] = 27951
$row7[0,2] = 27461
$ws.Range("AX7:AZ7").Value = $row7
$row8 = New-Object "object[,]" 1,3
$row8[0,0] = 0
$row8[0,1] = 0
$row8[0,2] = 0
$ws.Range("AX8:AZ8").Value = $row8
$row9 = New-Object "object[,]" 1,3
$row9[0,0] = 6434
$row9[0,1] = 4540
$row9[0,2] = 5204
$ws.Range("AX9:AZ9").Value = $row9
$row10 = New-Object "object[,]" 1,3
$row10[0,0] = 0
$row10[0,1] = 0
$row10[0,2] = 0
$ws.Range("AX10:AZ10").Value = $row10
$row11 = New-Object "object[,]" 1,3
$row11[0,0] = 9018
$row11[0,1] = 8360
$row11[0,2] = 9068
$ws.Range("AX11:AZ11").Value = $row11
$row12 = New-Object "object[,]" 1,3
$row12[0,0] = 33143
$row12[0,1] = 32248
$row12[0,2] = 32365
$ws.Range("AX12:AZ12").Value = $row12
$row13 = New-Object "object[,]" 1,3
$row13[0,0] = 0
$row13[0,1] = 0
$row13[0,2] = 0
$ws.Range("AX13:AZ13").Value = $row13
$row14 = New-Object "object[,]" 1,3
$row14[0,0] = 0
$row14[0,1] = 0
$row14[0,2] = 0
$ws.Range("AX14:AZ14").Value = $row14
$row15 = New-Object "object[,]" 1,3
$row15[0,0] = 10195
$row15[0,1] = 10482
$row15[0,2] = 10766
$ws.Range("AX15:AZ15").Value = $row15
$row16 = New-Object "object[,]" 1,3
$row16[0,0] = 13611
$row16[0,1] = 11948
$row16[0,2] = 10614
$ws.Range("AX16:AZ16").Value = $row16
$row17 = New-Object "object[,]" 1,3
$row17[0,0] = 0
$row17[0,1] = 0
$row17[0,2] = 0
$ws.Range("AX17:AZ17").Value = $row17
$row18 = New-Object "object[,]" 1,3
$row18[0,0] = 0
$row18[0,1] = 0
$row18[0,2] = 0
$ws.Range("AX18:AZ18").Value = $row18
$row19 = New-Object "object[,]" 1,3
$row19[0,0] = 3458
$row19[0,1] = 4117
$row19[0,2] = 5458
$ws.Range("AX19:AZ19").Value = $row19
$row20 = New-Object "object[,]" 1,3
$row20[0,0] = 0
$row20[0,1] = 0
$row20[0,2] = 0
$ws.Range("AX20:AZ20").Value = $row20
$row21 = New-Object "object[,]" 1,3
$row21[0,0] = 0
$row21[0,1] = 0
$row21[0,2] = 0
$ws.Range("AX21:AZ21").Value = $row21
$row22 = New-Object "object[,]" 1,3
$row22[0,0] = 20
$row22[0,1] = 20
$row22[0,2] = 20
$ws.Range("AX22:AZ22").Value = $row22
$row23 = New-Object "object[,]" 1,3
$row23[0,0] = 188746
$row23[0,1] = 192762
$row23[0,2] = 214578
$ws.Range("AX23:AZ23").Value = $row23
$row24 = New-Object "object[,]" 1,3
$row24[0,0] = 4439
$row24[0,1] = 4840
$row24[0,2] = 5023
$ws.Range("AX24:AZ24").Value = $row24
$row25 = New-Object "object[,]" 1,3
$row25[0,0] = 0
$row25[0,1] = 0
$row25[0,2] = 0
$ws.Range("AX25:AZ25").Value = $row25
$row26 = New-Object "object[,]" 1,3
$row26[0,0] = 349180.992
$row26[0,1] = 361632
$row26[0,2] = 357748
$ws.Range("AX26:AZ26").Value = $row26
$row27 = New-Object "object[,]" 1,3
$row27[0,0] = 83281
$row27[0,1] = 94181
$row27[0,2] = 82515
$ws.Range("AX27:AZ27").Value = $row27
$row28 = New-Object "object[,]" 1,3
$row28[0,0] = 5233
$row28[0,1] = 6395
$row28[0,2] = 9457
$ws.Range("AX28:AZ28").Value = $row28
$row29 = New-Object "object[,]" 1,3
$row29[0,0] = 6554
$row29[0,1] = 9534
$row29[0,2] = 13536
$ws.Range("AX29:AZ29").Value = $row29
$row30 = New-Object "object[,]" 1,3
$row30[0,0] = 1743
$row30[0,1] = 2664
$row30[0,2] = 5688
$ws.Range("AX30:AZ30").Value = $row30
$row31 = New-Object "object[,]" 1,3
$row31[0,0] = 7888
$row31[0,1] = 6867
$row31[0,2] = 5020
$ws.Range("AX31:AZ31").Value = $row31
$row32 = New-Object "object[,]" 1,3
$row32[0,0] = 0
$row32[0,1] = 0
$row32[0,2] = 0
$ws.Range("AX32:AZ32").Value = $row32
$row33 = New-Object "object[,]" 1,3
$row33[0,0] = 11672
$row33[0,1] = 11672
$row33[0,2] = 0
$ws.Range("AX33:AZ33").Value = $row33
$row34 = New-Object "object[,]" 1,3
$row34[0,0] = 50084
$row34[0,1] = 57049
$row34[0,2] = 48814
$ws.Range("AX34:AZ34").Value = $row34
$row35 = New-Object "object[,]" 1,3
$row35[0,0] = 107
$row35[0,1] = 0
$row35[0,2] = 0
$ws.Range("AX35:AZ35").Value = $row35
$row36 = New-Object "object[,]" 1,3
$row36[0,0] = 0
$row36[0,1] = 0
$row36[0,2] = 0
$ws.Range("AX36:AZ36").Value = $row36
$row37 = New-Object "object[,]" 1,3
$row37[0,0] = 75905
$row37[0,1] = 74396
$row37[0,2] = 78163
$ws.Range("AX37:AZ37").Value = $row37
$row38 = New-Object "object[,]" 1,3
$row38[0,0] = 69930
$row38[0,1] = 68747
$row38[0,2] = 71645
$ws.Range("AX38:AZ38").Value = $row38
$row39 = New-Object "object[,]" 1,3
$row39[0,0] = 0
$row39[0,1] = 0
$row39[0,2] = 0
$ws.Range("AX39:AZ39").Value = $row39
$row40 = New-Object "object[,]" 1,3
$row40[0,0] = 864
$row40[0,1] = 1002
$row40[0,2] = 622
$ws.Range("AX40:AZ40").Value = $row40
$row41 = New-Object "object[,]" 1,3
$row41[0,0] = 0
$row41[0,1] = 0
$row41[0,2] = 1225
$ws.Range("AX41:AZ41").Value = $row41
$row42 = New-Object "object[,]" 1,3
$row42[0,0] = 0
$row42[0,1] = 0
$row42[0,2] = 0
$ws.Range("AX42:AZ42").Value = $row42
$row43 = New-Object "object[,]" 1,3
$row43[0,0] = 5111
$row43[0,1] = 4647
$row43[0,2] = 4671
$ws.Range("AX43:AZ43").Value = $row43
$row44 = New-Object "object[,]" 1,3
$row44[0,0] = 0
$row44[0,1] = 0
$row44[0,2] = 0
$ws.Range("AX44:AZ44").Value = $row44
$row45 = New-Object "object[,]" 1,3
$row45[0,0] = 0
$row45[0,1] = 0
$row45[0,2] = 0
$ws.Range("AX45:AZ45").Value = $row45
$row46 = New-Object "object[,]" 1,3
$row46[0,0] = 0
$row46[0,1] = 0
$row46[0,2] = 0
$ws.Range("AX46:AZ46").Value = $row46
$row47 = New-Object "object[,]" 1,3
$row47[0,0] = 189995.008
$row47[0,1] = 193055.008
$row47[0,2] = 197070
$ws.Range("AX47:AZ47").Value = $row47
$row48 = New-Object "object[,]" 1,3
$row48[0,0] = 147000
$row48[0,1] = 147000
$row48[0,2] = 147000
$ws.Range("AX48:AZ48").Value = $row48
$row49 = New-Object "object[,]" 1,3
$row49[0,0] = 0
$row49[0,1] = 0
$row49[0,2] = 0
$ws.Range("AX49:AZ49").Value = $row49
$row50 = New-Object "object[,]" 1,3
$row50[0,0] = 0
$row50[0,1] = 0
$row50[0,2] = 0
$ws.Range("AX50:AZ50").Value = $row50
$row51 = New-Object "object[,]" 1,3
$row51[0,0] = 43284
$row51[0,1] = 43284
$row51[0,2] = 43284
$ws.Range("AX51:AZ51").Value = $row51
$row52 = New-Object "object[,]" 1,3
$row52[0,0] = 0
$row52[0,1] = 2803
$row52[0,2] = 6033
$ws.Range("AX52:AZ52").Value = $row52
$row53 = New-Object "object[,]" 1,3
$row53[0,0] = 0
$row53[0,1] = 0
$row53[0,2] = 0
$ws.Range("AX53:AZ53").Value = $row53
$row54 = New-Object "object[,]" 1,3
$row54[0,0] = -289
$row54[0,1] = -32
$row54[0,2] = 753
$ws.Range("AX54:AZ54").Value = $row54
$row55 = New-Object "object[,]" 1,3
$row55[0,0] = 0
$row55[0,1] = 0
$row55[0,2] = 0
$ws.Range("AX55:AZ55").Value = $row55
$row56 = New-Object "object[,]" 1,3
$row56[0,0] = 0
$row56[0,1] = 0
$row56[0,2] = 0
$ws.Range("AX56:AZ56").Value = $row56
$ws.Range("AX57:AZ57").Value = ""
$ws.Range("AX58:AZ58").Value = ""
$row59 = New-Object "object[,]" 1,3
$row59[0,0] = 58870.992
$row59[0,1] = 47093
$row59[0,2] = 56314
$ws.Range("AX59:AZ59").Value = $row59
$row60 = New-Object "object[,]" 1,3
$row60[0,0] = -37084.008
$row60[0,1] = -28437
$row60[0,2] = -34043
$ws.Range("AX60:AZ60").Value = $row60
$row61 = New-Object "object[,]" 1,3
$row61[0,0] = 21787
$row61[0,1] = 18656
$row61[0,2] = 22271
$ws.Range("AX61:AZ61").Value = $row61
$row62 = New-Object "object[,]" 1,3
$row62[0,0] = -16703
$row62[0,1] = -13657
$row62[0,2] = -15669
$ws.Range("AX62:AZ62").Value = $row62
$row63 = New-Object "object[,]" 1,3
$row63[0,0] = -5123
$row63[0,1] = -3961
$row63[0,2] = -4759
$ws.Range("AX63:AZ63").Value = $row63
$row64 = New-Object "object[,]" 1,3
$row64[0,0] = 0
$row64[0,1] = 0
$row64[0,2] = 0
$ws.Range("AX64:AZ64").Value = $row64
$row65 = New-Object "object[,]" 1,3
$row65[0,0] = 886
$row65[0,1] = 534
$row65[0,2] = 398
$ws.Range("AX65:AZ65").Value = $row65
$row66 = New-Object "object[,]" 1,3
$row66[0,0] = -70
$row66[0,1] = -573
$row66[0,2] = -5
$ws.Range("AX66:AZ66").Value = $row66
$row67 = New-Object "object[,]" 1,3
$row67[0,0] = 0
$row67[0,1] = 0
$row67[0,2] = 0
$ws.Range("AX67:AZ67").Value = $row67
$row68 = New-Object "object[,]" 1,3
$row68[0,0] = 1431
$row68[0,1] = 2383
$row68[0,2] = 2452
$ws.Range("AX68:AZ68").Value = $row68
$row69 = New-Object "object[,]" 1,3
$row69[0,0] = 3316
$row69[0,1] = 3424
$row69[0,2] = 4313
$ws.Range("AX69:AZ69").Value = $row69
$row70 = New-Object "object[,]" 1,3
$row70[0,0] = -1885
$row70[0,1] = -1041
$row70[0,2] = -1861
$ws.Range("AX70:AZ70").Value = $row70
$ws.Range("AX71:AZ71").Value = ""
$ws.Range("AX72:AZ72").Value = ""
$ws.Range("AX73:AZ73").Value = ""
$row74 = New-Object "object[,]" 1,3
$row74[0,0] = 2208
$row74[0,1] = 3382
$row74[0,2] = 4688
$ws.Range("AX74:AZ74").Value = $row74
$row75 = New-Object "object[,]" 1,3
$row75[0,0] = 2306
$row75[0,1] = -876
$row75[0,2] = -1335
$ws.Range("AX75:AZ75").Value = $row75
$row76 = New-Object "object[,]" 1,3
$row76[0,0] = 1995
$row76[0,1] = 297
$row76[0,2] = -123
$ws.Range("AX76:AZ76").Value = $row76
$ws.Range("AX77:AZ77").Value = ""
$ws.Range("AX78:AZ78").Value = ""
$row79 = New-Object "object[,]" 1,3
$row79[0,0] = 0
$row79[0,1] = 0
$row79[0,2] = 0
$ws.Range("AX79:AZ79").Value = $row79
$row80 = New-Object "object[,]" 1,3
$row80[0,0] = 6509
$row80[0,1] = 2803
$row80[0,2] = 3230
$ws.Range("AX80:AZ80").Value = $row80
